$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").Borders.Weight = 2

$iValues = @(7,6,6,1,7,5,6,7,1,1,1,1)
$jValues = @(8,6,6,2,8,7,7,8,2,3,1,1)

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
